$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.321.38'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '1.841.88'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('D4').Value = "'0.9990"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'238.75"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('D6').Value = "'0.6300"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'0.07521"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.93%  '
$ws.Range('D9').Value = "'0.2941"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('D10').Value = "'24.44"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').Value = "'0.07692"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = '1.833.45'
$ws.Range('E12').Value = '  -7.64%  '
$ws.Range('D14').Value = "'0.6779"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = "'0.00001045"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.47%  '
$ws.Range('D16').Value = "'82.96"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = "'6.132"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '29.356.84'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = "'228.14"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.48%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = "'12.42"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.80%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = "'0.9997"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = "'7.412"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.57%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = "'1.000"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = "'156.60"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.46%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').Value = "'0.1389"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = "'8.342"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.43%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'17.59"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = "'1.454"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.38%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'1.271"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = "'0.05631"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.16%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'4.102"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'4.017"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.24%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').Value = "'1.830"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.27%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = "'1.154"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = "'0.7085"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.590"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.244.86'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'0.01809"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = "'2.758"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.16%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'6.238"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = "'0.9021"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = "'0.9992"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = "'101.85"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = "'65.58"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.52%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = "'0.00000000120"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.37%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = "'7.095"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.69%  '
$ws.Range('B47').Value = 'TheSandbox'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D47').Value = "'0.3990"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'8.923"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.97%  '
$ws.Range('D49').Value = "'1.670"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.50%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = "'0.1120"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.05712"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.65%  '
